$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: snapshot the two existing style templates into scratch rows
# before we overwrite rows 2-4 (pattern P1: F=fail/red,G-J=pass/green)
$ws.Range("A2:J2").Copy()
$ws.Range("A100:J100").PasteSpecial(-4122)
# pattern P2: F=pass/green, G-J=fail/red
$ws.Range("A3:J3").Copy()
$ws.Range("A101:J101").PasteSpecial(-4122)

# --- Step 2: apply format templates to every data row (2-12) first
$ws.Range("A100:J100").Copy()
$ws.Range("A2:J2").PasteSpecial(-4122)
$ws.Range("A100:J100").Copy()
$ws.Range("A3:J3").PasteSpecial(-4122)
$ws.Range("A100:J100").Copy()
$ws.Range("A4:J4").PasteSpecial(-4122)
$ws.Range("A101:J101").Copy()
$ws.Range("A5:J5").PasteSpecial(-4122)
$ws.Range("A101:J101").Copy()
$ws.Range("A6:J6").PasteSpecial(-4122)
$ws.Range("A101:J101").Copy()
$ws.Range("A7:J7").PasteSpecial(-4122)
$ws.Range("A100:J100").Copy()
$ws.Range("A8:J8").PasteSpecial(-4122)
# row 9 has a unique pattern ('2', '2', '2', '2', '2', '4', '3', '4', '3', '3'), handled per-cell below
$ws.Range("A101:J101").Copy()
$ws.Range("A10:J10").PasteSpecial(-4122)
$ws.Range("A101:J101").Copy()
$ws.Range("A11:J11").PasteSpecial(-4122)
$ws.Range("A101:J101").Copy()
$ws.Range("A12:J12").PasteSpecial(-4122)

# row 9 unique pattern differs from template P2 only at H9 (style 4 instead of 3)
$ws.Range("A101:J101").Copy()
$ws.Range("A9:J9").PasteSpecial(-4122)
$ws.Range("G2").Copy()
$ws.Range("H9").PasteSpecial(-4122)

# --- Step 3: clear the scratch template rows
$ws.Range("A100:J101").Clear()

# --- Step 4: write values
$ws.Range("A1").Value = "file_name"
$ws.Range("B1").Value = "file_executions"
$ws.Range("C1").Value = "sheet_name"
$ws.Range("D1").Value = "executed_sheets"
$ws.Range("E1").Value = "total_rows"
$ws.Range("F1").Value = "pass"
$ws.Range("G1").Value = "number_fail"
$ws.Range("H1").Value = "key_fail"
$ws.Range("I1").Value = "sum_value_differences"
$ws.Range("J1").Value = "max_difference"
$ws.Range("A2").Value = "MOCK_DATA_old - Copy here is very big name of the file, what can extend the total table.xlsx, file missed: NOT EXECUTED"
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = "NONE"
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("A3").Value = "MOCK_DATA_old.xlsx, file missed: NOT EXECUTED"
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = "NONE"
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("A4").Value = "MOCK_DATA_old_csv.csv, file missed: NOT EXECUTED"
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = "NONE"
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("A5").Value = "banana_quality_dataset.xlsx"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = "banana_quality_dataset.csv"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 976
$ws.Range("F5").Value = 967
$ws.Range("G5").Value = 5
$ws.Range("H5").Value = 4
$ws.Range("I5").Value = -2343.4462
$ws.Range("J5").Value = 2210
$ws.Range("A6").Value = "MOCK_DATA - csv.csv"
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = "sheet"
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 152
$ws.Range("F6").Value = 144
$ws.Range("G6").Value = 3
$ws.Range("H6").Value = 5
$ws.Range("I6").Value = -31671.00874
$ws.Range("J6").Value = 31190
$ws.Range("A7").Value = "Test_500_without big number.xlsx"
$ws.Range("B7").Value = 3
$ws.Range("C7").Value = "MOCK_DATA (1).csv"
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 504
$ws.Range("F7").Value = 493
$ws.Range("G7").Value = 3
$ws.Range("H7").Value = 8
$ws.Range("I7").Value = 5088
$ws.Range("J7").Value = 9900
$ws.Range("A8").Value = "MOCK_DATA.xlsx"
$ws.Range("B8").Value = 4
$ws.Range("C8").Value = "['extra_sheet'] missing: NOT EXECUTED"
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("A9").Value = "MOCK_DATA.xlsx"
$ws.Range("B9").Value = 4
$ws.Range("C9").Value = "some_data"
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 49
$ws.Range("F9").Value = 46
$ws.Range("G9").Value = 3
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 45.95812
$ws.Range("J9").Value = 46
$ws.Range("A10").Value = "MOCK_DATA.xlsx"
$ws.Range("B10").Value = 4
$ws.Range("C10").Value = "mock_data"
$ws.Range("D10").Value = 2
$ws.Range("E10").Value = 152
$ws.Range("F10").Value = 144
$ws.Range("G10").Value = 4
$ws.Range("H10").Value = 4
$ws.Range("I10").Value = 114.00421
$ws.Range("J10").Value = 70
$ws.Range("A11").Value = "Test_500.xlsx"
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = "MOCK_DATA (1).csv"
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 502
$ws.Range("F11").Value = 494
$ws.Range("G11").Value = 4
$ws.Range("H11").Value = 4
$ws.Range("I11").Value = 665488751
$ws.Range("J11").Value = 730893583
$ws.Range("A12").Value = "TOTAL"
$ws.Range("B12").Value = 5
$ws.Range("C12").Value = ""
$ws.Range("D12").Value = 6
$ws.Range("E12").Value = 2335
$ws.Range("F12").Value = 2288
$ws.Range("G12").Value = 22
$ws.Range("H12").Value = 25
$ws.Range("I12").Value = 665459984.50739
$ws.Range("J12").Value = 730893583

# --- Step 5: column A width -> stored width 121
$ws.Columns.Item(1).ColumnWidth = 120.16666666666666
